# Update generalized_cf sheet: refresh simulated values for rows 2-4 and
# append a new row 5 (bucket index 3), matching the upstream re-run of the
# algo. Column A holds the bucket index (0,1,2,3,...); columns B:BA hold the
# 50 per-trial values; AZ/BA (already present as headers "mean"/"std") are
# laid out one column before the 50-value block ends, matching the existing
# sheet layout exactly as in the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for existing rows 2, 3 and 4 (columns B..BA -> 2..53)
$row2 = @(0.00639817195087118,0.007084499800034279,0.006871278057718736,0.007956952315530367,0.006976212259835315,0.006921404873584258,0.007432818753573471,0.007097882083571838,0.006173545215502458,0.006346846589284693,0.006807780320366132,0.007328523989465247,0.007267525035765379,0.007786110952081067,0.006689536878216123,0.007024557395773843,0.007423904974016332,0.006518010291595197,0.005828571428571429,0.007024156244646223,0.006466380543633763,0.006586106179485711,0.007196710075394105,0.006870491240123669,0.006573306659045442,0.007438347542484408,0.006745941001600732,0.007494279176201373,0.007440050363417845,0.007098288396588242,0.007093009953094612,0.005998286203941731,0.006578194714563551,0.006863025450386045,0.007605649911362727,0.007333982696384576,0.007662835249042145,0.006358480838632067,0.006064420161336461,0.0073735352957988,0.007090170964606324,0.006292546193009553,0.0072,0.006858710562414266,0.006755596267246808,0.007037015847588535,0.006861455772199669,0.007897899616551251,0.00760826039700246,0.006753276483717736,0.006963090863337165,0.0004919310072975257)
$row3 = @(0.1845140032948929,0.2087542087542087,0.1935483870967742,0.2274959083469722,0.2023217247097844,0.2023411371237458,0.2199661590524535,0.2059800664451827,0.1821247892074199,0.1837748344370861,0.1980033277870216,0.2088091353996737,0.2123745819397993,0.2233169129720854,0.1875,0.2043189368770764,0.2155887230514096,0.1856677524429967,0.1746575342465753,0.2067226890756302,0.1928327645051195,0.1959114139693356,0.2058823529411765,0.2051282051282051,0.1975945017182131,0.2152317880794702,0.1976549413735343,0.2183333333333333,0.2233676975945017,0.2087542087542087,0.2003231017770598,0.175,0.195578231292517,0.1973684210526316,0.2205638474295191,0.2140468227424749,0.2126984126984127,0.1816693944353519,0.1729200652528548,0.2135761589403974,0.209106239460371,0.1824212271973466,0.2052117263843648,0.2006688963210702,0.193126022913257,0.2,0.1990049751243781,0.2262295081967213,0.2224080267558528,0.200339558573854,0.2022146531241265,0.01390095684249412)
$row4 = @(0.01236749116607774,0.01370392882798254,0.013271400132714,0.01537610619469026,0.01348736941020397,0.01338495575221239,0.01437973563409104,0.01372288623284639,0.01194227898490629,0.01226993865030675,0.01316298877274487,0.01416007522539964,0.01405411387152105,0.0150475768975437,0.01291818482941371,0.0135821554770318,0.01435353869934857,0.01259390190013257,0.0112806901128069,0.01358665635700873,0.01251314988095897,0.01274379432624114,0.01390728476821192,0.01329566229017783,0.01272335011340377,0.01437973563409104,0.01304660290784455,0.01449115044247788,0.01440044309055663,0.01372972374467143,0.01370090050273466,0.0115990057995029,0.01272827891532928,0.01326479854087216,0.01470425649530127,0.01418203977618969,0.01479273610421151,0.01228691609475315,0.01171788635861154,0.01425493121166915,0.01371529698042252,0.0121654501216545,0.01391189135475323,0.01326406543605615,0.01305454143157429,0.01359566707195755,0.01326553172673005,0.01526295415583697,0.01471320316389181,0.01306610563614218,0.01346246654411625,0.0009492699078558605)

# Brand-new row 5 (bucket index 3), same column layout as rows 2-4
$row5 = @(0.7860302677532014,0.7927823050058207,0.7850989522700815,0.7862630966239814,0.7913853317811409,0.7825378346915017,0.7762514551804424,0.79953434225844,0.7818393480791618,0.7783469150174621,0.7913853317811409,0.7925494761350408,0.7871944121071013,0.790221187427241,0.7844004656577416,0.780442374854482,0.7864959254947613,0.7885913853317812,0.7881257275902211,0.7867287543655413,0.7918509895227008,0.7904540162980209,0.7846332945285215,0.7809080325960419,0.7885913853317812,0.790221187427241,0.7834691501746216,0.7811408614668219,0.7878928987194412,0.7883585564610012,0.7883585564610012,0.7878928987194412,0.7874272409778813,0.7813736903376018,0.7750873108265425,0.7790454016298021,0.7841676367869616,0.8013969732246798,0.7813736903376018,0.7834691501746216,0.780442374854482,0.789988358556461,0.7869615832363213,0.7846332945285215,0.7867287543655413,0.7881257275902211,0.7953434225844005,0.7823050058207218,0.7848661233993015,0.7788125727590222,0.7862305005820722,0.005288242300914685)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
}
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
}
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 2).Value = $row4[$i]
}

# Row 5: bucket index in column A, then the 50 trial values + mean/std
$ws.Cells.Item(5, 1).Value = 3
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $row5[$i]
}

# Give the new index cell (A5) the same bold/bordered style as A2:A4
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
